$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end, after "Better split" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Cherry Picked Split"

# --- Populate data rows first (2-68) so the shared-string table picks up
#     the image file names before the two brand-new header labels ---
$newSheet.Range("A2").Value = "N10_L.png"
$newSheet.Range("B2").Value = 0
$newSheet.Range("C2").Value = 0
$newSheet.Range("D2").Value = 0
$newSheet.Range("E2").Value = 0.38127225621974398
$newSheet.Range("F2").Value = 0.462857142
$newSheet.Range("G2").Value = 0
$newSheet.Range("A3").Value = "N10_R.png"
$newSheet.Range("B3").Value = 0
$newSheet.Range("C3").Value = 0
$newSheet.Range("D3").Value = 0
$newSheet.Range("E3").Value = 0.35214388457679602
$newSheet.Range("F3").Value = 0.305714285999999
$newSheet.Range("G3").Value = 0
$newSheet.Range("A4").Value = "N11_L.png"
$newSheet.Range("B4").Value = 0
$newSheet.Range("C4").Value = 0
$newSheet.Range("D4").Value = 0
$newSheet.Range("E4").Value = 0.42568027393660302
$newSheet.Range("F4").Value = 0.54666666600000002
$newSheet.Range("G4").Value = 0
$newSheet.Range("A5").Value = "N11_R.png"
$newSheet.Range("B5").Value = 0
$newSheet.Range("C5").Value = 1
$newSheet.Range("D5").Value = 0
$newSheet.Range("E5").Value = 0.48328903382826499
$newSheet.Range("F5").Value = 0.40571428599999898
$newSheet.Range("G5").Value = 1
$newSheet.Range("A6").Value = "N12_L.png"
$newSheet.Range("B6").Value = 0
$newSheet.Range("C6").Value = 1
$newSheet.Range("D6").Value = 0
$newSheet.Range("E6").Value = 0.56278733059599895
$newSheet.Range("F6").Value = 0.182857142
$newSheet.Range("G6").Value = 1
$newSheet.Range("A7").Value = "N13_L.png"
$newSheet.Range("B7").Value = 0
$newSheet.Range("C7").Value = 0
$newSheet.Range("D7").Value = 0
$newSheet.Range("E7").Value = 0.46802736164238101
$newSheet.Range("F7").Value = 0.29857142799999897
$newSheet.Range("G7").Value = 0
$newSheet.Range("A8").Value = "N13_R.png"
$newSheet.Range("B8").Value = 0
$newSheet.Range("C8").Value = 0
$newSheet.Range("D8").Value = 0
$newSheet.Range("E8").Value = 0.278069908418352
$newSheet.Range("F8").Value = 0.55466666599999903
$newSheet.Range("G8").Value = 0
$newSheet.Range("A9").Value = "N14_R.png"
$newSheet.Range("B9").Value = 0
$newSheet.Range("C9").Value = 0
$newSheet.Range("D9").Value = 0
$newSheet.Range("E9").Value = 0.126904712172594
$newSheet.Range("F9").Value = 0.342857142
$newSheet.Range("G9").Value = 0
$newSheet.Range("A10").Value = "N16_L.png"
$newSheet.Range("B10").Value = 0
$newSheet.Range("C10").Value = 1
$newSheet.Range("D10").Value = 0
$newSheet.Range("E10").Value = 0.34085778617495299
$newSheet.Range("F10").Value = 0.0414285719999999
$newSheet.Range("G10").Value = 1
$newSheet.Range("A11").Value = "N16_R.png"
$newSheet.Range("B11").Value = 0
$newSheet.Range("C11").Value = 0
$newSheet.Range("D11").Value = 0
$newSheet.Range("E11").Value = 0.122874550414414
$newSheet.Range("F11").Value = 0.47538461599999998
$newSheet.Range("G11").Value = 0
$newSheet.Range("A12").Value = "N17_L.png"
$newSheet.Range("B12").Value = 0
$newSheet.Range("C12").Value = 0
$newSheet.Range("D12").Value = 0
$newSheet.Range("E12").Value = 0.122979515794961
$newSheet.Range("F12").Value = 0.54933333399999995
$newSheet.Range("G12").Value = 0
$newSheet.Range("A13").Value = "N17_R.png"
$newSheet.Range("B13").Value = 0
$newSheet.Range("C13").Value = 0
$newSheet.Range("D13").Value = 1
$newSheet.Range("E13").Value = 0.45815132048564999
$newSheet.Range("F13").Value = 0.084000000000000005
$newSheet.Range("G13").Value = 0
$newSheet.Range("A14").Value = "N18_L.png"
$newSheet.Range("B14").Value = 0
$newSheet.Range("C14").Value = 0
$newSheet.Range("D14").Value = 0
$newSheet.Range("E14").Value = 0.026363433362033301
$newSheet.Range("F14").Value = 0.42399999999999999
$newSheet.Range("G14").Value = 0
$newSheet.Range("A15").Value = "N18_R.png"
$newSheet.Range("B15").Value = 0
$newSheet.Range("C15").Value = 0
$newSheet.Range("D15").Value = 0
$newSheet.Range("E15").Value = 0.67731031878068804
$newSheet.Range("F15").Value = 0.347142858
$newSheet.Range("G15").Value = 0
$newSheet.Range("A16").Value = "N19_L.png"
$newSheet.Range("B16").Value = 0
$newSheet.Range("C16").Value = 0
$newSheet.Range("D16").Value = 0
$newSheet.Range("E16").Value = 1.93534490534277
$newSheet.Range("F16").Value = 0.61833333399999901
$newSheet.Range("G16").Value = 0
$newSheet.Range("A17").Value = "N19_R.png"
$newSheet.Range("B17").Value = 0
$newSheet.Range("C17").Value = 0
$newSheet.Range("D17").Value = 0
$newSheet.Range("E17").Value = 1.34905680728227
$newSheet.Range("F17").Value = 0.34153846199999999
$newSheet.Range("G17").Value = 0
$newSheet.Range("A18").Value = "N1_L.png"
$newSheet.Range("B18").Value = 0
$newSheet.Range("C18").Value = 0
$newSheet.Range("D18").Value = 0
$newSheet.Range("E18").Value = 1.39494353385168
$newSheet.Range("F18").Value = 0.37066666599999898
$newSheet.Range("G18").Value = 0
$newSheet.Range("A19").Value = "N1_R.png"
$newSheet.Range("B19").Value = 0
$newSheet.Range("C19").Value = 0
$newSheet.Range("D19").Value = 0
$newSheet.Range("E19").Value = 0.74689730504989105
$newSheet.Range("F19").Value = 0.107999999999999
$newSheet.Range("G19").Value = 0
$newSheet.Range("A20").Value = "N23_L.png"
$newSheet.Range("B20").Value = 0
$newSheet.Range("C20").Value = 0
$newSheet.Range("D20").Value = 0
$newSheet.Range("E20").Value = 1.0038649892967999
$newSheet.Range("F20").Value = 0.62714285800000003
$newSheet.Range("G20").Value = 0
$newSheet.Range("A21").Value = "N23_R.png"
$newSheet.Range("B21").Value = 0
$newSheet.Range("C21").Value = 1
$newSheet.Range("D21").Value = 0
$newSheet.Range("E21").Value = 0.50606762399137595
$newSheet.Range("F21").Value = 0.55333333399999896
$newSheet.Range("G21").Value = 0
$newSheet.Range("A22").Value = "N2_L.png"
$newSheet.Range("B22").Value = 0
$newSheet.Range("C22").Value = 1
$newSheet.Range("D22").Value = 0
$newSheet.Range("E22").Value = 0.345858421257078
$newSheet.Range("F22").Value = 0.37230769199999902
$newSheet.Range("G22").Value = 0
$newSheet.Range("A23").Value = "N2_R.png"
$newSheet.Range("B23").Value = 0
$newSheet.Range("C23").Value = 1
$newSheet.Range("D23").Value = 0
$newSheet.Range("E23").Value = 0.418320598708599
$newSheet.Range("F23").Value = 0.104
$newSheet.Range("G23").Value = 1
$newSheet.Range("A24").Value = "N3_L.png"
$newSheet.Range("B24").Value = 0
$newSheet.Range("C24").Value = 1
$newSheet.Range("D24").Value = 0
$newSheet.Range("E24").Value = 0.26355140558377099
$newSheet.Range("F24").Value = 0.24153846200000001
$newSheet.Range("G24").Value = 1
$newSheet.Range("A25").Value = "N3_R.png"
$newSheet.Range("B25").Value = 0
$newSheet.Range("C25").Value = 1
$newSheet.Range("D25").Value = 0
$newSheet.Range("E25").Value = 0.489928809122961
$newSheet.Range("F25").Value = 0.31733333399999902
$newSheet.Range("G25").Value = 1
$newSheet.Range("A26").Value = "N4_L.png"
$newSheet.Range("B26").Value = 0
$newSheet.Range("C26").Value = 0
$newSheet.Range("D26").Value = 0
$newSheet.Range("E26").Value = 1.61310979439297
$newSheet.Range("F26").Value = 0.59599999999999997
$newSheet.Range("G26").Value = 0
$newSheet.Range("A27").Value = "N5_L.png"
$newSheet.Range("B27").Value = 0
$newSheet.Range("C27").Value = 0
$newSheet.Range("D27").Value = 0
$newSheet.Range("E27").Value = 0.071176302689889101
$newSheet.Range("F27").Value = 0.34399999999999997
$newSheet.Range("G27").Value = 0
$newSheet.Range("A28").Value = "N5_R.png"
$newSheet.Range("B28").Value = 0
$newSheet.Range("C28").Value = 0
$newSheet.Range("D28").Value = 0
$newSheet.Range("E28").Value = 0.78841923541751202
$newSheet.Range("F28").Value = 0.54533333399999995
$newSheet.Range("G28").Value = 0
$newSheet.Range("A29").Value = "N6_L.png"
$newSheet.Range("B29").Value = 0
$newSheet.Range("C29").Value = 0
$newSheet.Range("D29").Value = 0
$newSheet.Range("E29").Value = 0.0144206486581063
$newSheet.Range("F29").Value = 0.62142857200000001
$newSheet.Range("G29").Value = 0
$newSheet.Range("A30").Value = "N7_L.png"
$newSheet.Range("B30").Value = 0
$newSheet.Range("C30").Value = 1
$newSheet.Range("D30").Value = 0
$newSheet.Range("E30").Value = 0.945894960519573
$newSheet.Range("F30").Value = 0.49333333399999901
$newSheet.Range("G30").Value = 1
$newSheet.Range("A31").Value = "N7_R.png"
$newSheet.Range("B31").Value = 0
$newSheet.Range("C31").Value = 1
$newSheet.Range("D31").Value = 0
$newSheet.Range("E31").Value = 0.0014421543312912301
$newSheet.Range("F31").Value = 0.133333334
$newSheet.Range("G31").Value = 0
$newSheet.Range("A32").Value = "N8_R.png"
$newSheet.Range("B32").Value = 0
$newSheet.Range("C32").Value = 0
$newSheet.Range("D32").Value = 0
$newSheet.Range("E32").Value = 0.276114718525655
$newSheet.Range("F32").Value = 0.55733333399999996
$newSheet.Range("G32").Value = 0
$newSheet.Range("A33").Value = "N9_L.png"
$newSheet.Range("B33").Value = 0
$newSheet.Range("C33").Value = 1
$newSheet.Range("D33").Value = 0
$newSheet.Range("E33").Value = 0.014005968750279299
$newSheet.Range("F33").Value = 0.42307692400000002
$newSheet.Range("G33").Value = 0
$newSheet.Range("A34").Value = "N9_R.png"
$newSheet.Range("B34").Value = 0
$newSheet.Range("C34").Value = 0
$newSheet.Range("D34").Value = 1
$newSheet.Range("E34").Value = 0.83427774631806295
$newSheet.Range("F34").Value = 0.15428571399999899
$newSheet.Range("G34").Value = 0
$newSheet.Range("A35").Value = "AD13_L.png"
$newSheet.Range("B35").Value = 1
$newSheet.Range("C35").Value = 0
$newSheet.Range("D35").Value = 1
$newSheet.Range("E35").Value = 1.17545420169445
$newSheet.Range("F35").Value = 0.59285714199999995
$newSheet.Range("G35").Value = 0
$newSheet.Range("A36").Value = "AD15_R.png"
$newSheet.Range("B36").Value = 1
$newSheet.Range("C36").Value = 1
$newSheet.Range("D36").Value = 1
$newSheet.Range("E36").Value = 0.13469359824778501
$newSheet.Range("F36").Value = 0.56666666600000004
$newSheet.Range("G36").Value = 1
$newSheet.Range("A37").Value = "AD16_R.png"
$newSheet.Range("B37").Value = 1
$newSheet.Range("C37").Value = 1
$newSheet.Range("D37").Value = 0
$newSheet.Range("E37").Value = 0.35679965585986001
$newSheet.Range("F37").Value = 0.29428571399999998
$newSheet.Range("G37").Value = 1
$newSheet.Range("A38").Value = "AD17_L.png"
$newSheet.Range("B38").Value = 1
$newSheet.Range("C38").Value = 0
$newSheet.Range("D38").Value = 1
$newSheet.Range("E38").Value = 1.2527240589111399
$newSheet.Range("F38").Value = 0.14000000000000001
$newSheet.Range("G38").Value = 0
$newSheet.Range("A39").Value = "AD18_R.png"
$newSheet.Range("B39").Value = 1
$newSheet.Range("C39").Value = 1
$newSheet.Range("D39").Value = 1
$newSheet.Range("E39").Value = 0.83770278143471
$newSheet.Range("F39").Value = 0.18266666599999901
$newSheet.Range("G39").Value = 1
$newSheet.Range("A40").Value = "AD19_L.png"
$newSheet.Range("B40").Value = 1
$newSheet.Range("C40").Value = 1
$newSheet.Range("D40").Value = 1
$newSheet.Range("E40").Value = 0.132856754599637
$newSheet.Range("F40").Value = 0.029333333999999898
$newSheet.Range("G40").Value = 1
$newSheet.Range("A41").Value = "AD1_L.png"
$newSheet.Range("B41").Value = 1
$newSheet.Range("C41").Value = 0
$newSheet.Range("D41").Value = 1
$newSheet.Range("E41").Value = 0.039194285037401501
$newSheet.Range("F41").Value = 0.62428571399999999
$newSheet.Range("G41").Value = 1
$newSheet.Range("A42").Value = "AD20_L.png"
$newSheet.Range("B42").Value = 1
$newSheet.Range("C42").Value = 1
$newSheet.Range("D42").Value = 1
$newSheet.Range("E42").Value = 0.90424178668393196
$newSheet.Range("F42").Value = 0.42249999999999999
$newSheet.Range("G42").Value = 1
$newSheet.Range("A43").Value = "AD22_L.png"
$newSheet.Range("B43").Value = 1
$newSheet.Range("C43").Value = 0
$newSheet.Range("D43").Value = 0
$newSheet.Range("E43").Value = 0.434172685427821
$newSheet.Range("F43").Value = 0.16666666599999999
$newSheet.Range("G43").Value = 0
$newSheet.Range("A44").Value = "AD23_L.png"
$newSheet.Range("B44").Value = 1
$newSheet.Range("C44").Value = 1
$newSheet.Range("D44").Value = 1
$newSheet.Range("E44").Value = 0.29161066787554502
$newSheet.Range("F44").Value = 0.71599999999999997
$newSheet.Range("G44").Value = 1
$newSheet.Range("A45").Value = "AD25_L.png"
$newSheet.Range("B45").Value = 1
$newSheet.Range("C45").Value = 1
$newSheet.Range("D45").Value = 0
$newSheet.Range("E45").Value = 0.47851370586849901
$newSheet.Range("F45").Value = 0.0046153839999999403
$newSheet.Range("G45").Value = 1
$newSheet.Range("A46").Value = "AD27_R.png"
$newSheet.Range("B46").Value = 1
$newSheet.Range("C46").Value = 0
$newSheet.Range("D46").Value = 1
$newSheet.Range("E46").Value = 0.41839552317304002
$newSheet.Range("F46").Value = 0.08
$newSheet.Range("G46").Value = 0
$newSheet.Range("A47").Value = "AD29_R.png"
$newSheet.Range("B47").Value = 1
$newSheet.Range("C47").Value = 1
$newSheet.Range("D47").Value = 1
$newSheet.Range("E47").Value = 0.46462417448906701
$newSheet.Range("F47").Value = 0.342857142
$newSheet.Range("G47").Value = 1
$newSheet.Range("A48").Value = "AD31_R.png"
$newSheet.Range("B48").Value = 1
$newSheet.Range("C48").Value = 0
$newSheet.Range("D48").Value = 1
$newSheet.Range("E48").Value = 0.102525642531787
$newSheet.Range("F48").Value = 0.51200000000000001
$newSheet.Range("G48").Value = 1
$newSheet.Range("A49").Value = "AD33_L.png"
$newSheet.Range("B49").Value = 1
$newSheet.Range("C49").Value = 1
$newSheet.Range("D49").Value = 1
$newSheet.Range("E49").Value = 0.57153831968542002
$newSheet.Range("F49").Value = 0.57428571399999995
$newSheet.Range("G49").Value = 1
$newSheet.Range("A50").Value = "AD34_R.png"
$newSheet.Range("B50").Value = 1
$newSheet.Range("C50").Value = 1
$newSheet.Range("D50").Value = 1
$newSheet.Range("E50").Value = 0.12919397232232599
$newSheet.Range("F50").Value = 0.14000000000000001
$newSheet.Range("G50").Value = 1
$newSheet.Range("A51").Value = "AD36_R.png"
$newSheet.Range("B51").Value = 1
$newSheet.Range("C51").Value = 1
$newSheet.Range("D51").Value = 1
$newSheet.Range("E51").Value = 0.31799572002711102
$newSheet.Range("F51").Value = 0.42266666600000002
$newSheet.Range("G51").Value = 1
$newSheet.Range("A52").Value = "AD3_L.png"
$newSheet.Range("B52").Value = 1
$newSheet.Range("C52").Value = 0
$newSheet.Range("D52").Value = 1
$newSheet.Range("E52").Value = 0.59193187084064403
$newSheet.Range("F52").Value = 0.30249999999999999
$newSheet.Range("G52").Value = 0
$newSheet.Range("A53").Value = "AD48_R.png"
$newSheet.Range("B53").Value = 1
$newSheet.Range("C53").Value = 1
$newSheet.Range("D53").Value = 0
$newSheet.Range("E53").Value = 0.306067726225213
$newSheet.Range("F53").Value = 0.06
$newSheet.Range("G53").Value = 1
$newSheet.Range("A54").Value = "AD4_R.png"
$newSheet.Range("B54").Value = 1
$newSheet.Range("C54").Value = 1
$newSheet.Range("D54").Value = 1
$newSheet.Range("E54").Value = 0.47309183145650802
$newSheet.Range("F54").Value = 0.27
$newSheet.Range("G54").Value = 1
$newSheet.Range("A55").Value = "AD15_L.png"
$newSheet.Range("B55").Value = "CONTRALATERAL"
$newSheet.Range("C55").Value = 1
$newSheet.Range("D55").Value = 1
$newSheet.Range("E55").Value = 0.21275148197621599
$newSheet.Range("F55").Value = 0.055714286000000002
$newSheet.Range("A56").Value = "AD22_R.png"
$newSheet.Range("B56").Value = "CONTRALATERAL"
$newSheet.Range("C56").Value = 1
$newSheet.Range("D56").Value = 0
$newSheet.Range("E56").Value = 0.47628868300782101
$newSheet.Range("F56").Value = 0.34615384599999999
$newSheet.Range("A57").Value = "AD27_L.png"
$newSheet.Range("B57").Value = "CONTRALATERAL"
$newSheet.Range("C57").Value = 1
$newSheet.Range("D57").Value = 0
$newSheet.Range("E57").Value = 0.23080995081538799
$newSheet.Range("F57").Value = 0.28533333399999999
$newSheet.Range("A58").Value = "AD29_L.png"
$newSheet.Range("B58").Value = "CONTRALATERAL"
$newSheet.Range("C58").Value = 1
$newSheet.Range("D58").Value = 0
$newSheet.Range("E58").Value = 0.65017711595093497
$newSheet.Range("F58").Value = 0.20714285799999901
$newSheet.Range("A59").Value = "AD2_R.png"
$newSheet.Range("B59").Value = "CONTRALATERAL"
$newSheet.Range("C59").Value = 1
$newSheet.Range("D59").Value = 0
$newSheet.Range("E59").Value = 1.08540133661293
$newSheet.Range("F59").Value = 0.57571428599999896
$newSheet.Range("A60").Value = "AD32_R.png"
$newSheet.Range("B60").Value = "CONTRALATERAL"
$newSheet.Range("C60").Value = 0
$newSheet.Range("D60").Value = 0
$newSheet.Range("E60").Value = 0.133919099500029
$newSheet.Range("F60").Value = 0.010666666
$newSheet.Range("A61").Value = "AD35_L.png"
$newSheet.Range("B61").Value = "CONTRALATERAL"
$newSheet.Range("C61").Value = 1
$newSheet.Range("D61").Value = 1
$newSheet.Range("E61").Value = 0.45963143998139599
$newSheet.Range("F61").Value = 0.17733333400000001
$newSheet.Range("A62").Value = "AD36_L.png"
$newSheet.Range("B62").Value = "CONTRALATERAL"
$newSheet.Range("C62").Value = 0
$newSheet.Range("D62").Value = 0
$newSheet.Range("E62").Value = 0.57380544289633895
$newSheet.Range("F62").Value = 0.4
$newSheet.Range("A63").Value = "AD38_R.png"
$newSheet.Range("B63").Value = "CONTRALATERAL"
$newSheet.Range("C63").Value = 1
$newSheet.Range("D63").Value = 0
$newSheet.Range("E63").Value = 1.02323805176822
$newSheet.Range("F63").Value = 0.31066666599999898
$newSheet.Range("A64").Value = "AD3_R.png"
$newSheet.Range("B64").Value = "CONTRALATERAL"
$newSheet.Range("C64").Value = 0
$newSheet.Range("D64").Value = 1
$newSheet.Range("E64").Value = 0.65711276624159198
$newSheet.Range("F64").Value = 0.270666666
$newSheet.Range("A65").Value = "AD47_L.png"
$newSheet.Range("B65").Value = "CONTRALATERAL"
$newSheet.Range("C65").Value = 0
$newSheet.Range("D65").Value = 0
$newSheet.Range("E65").Value = 0.097585339335031607
$newSheet.Range("F65").Value = 0.054285713999999902
$newSheet.Range("A66").Value = "AD6_L.png"
$newSheet.Range("B66").Value = "CONTRALATERAL"
$newSheet.Range("C66").Value = 0
$newSheet.Range("D66").Value = 0
$newSheet.Range("E66").Value = 0.33678689329197098
$newSheet.Range("F66").Value = 0.21230769199999899
$newSheet.Range("A67").Value = "AD7_L.png"
$newSheet.Range("B67").Value = "CONTRALATERAL"
$newSheet.Range("C67").Value = 1
$newSheet.Range("D67").Value = 1
$newSheet.Range("E67").Value = 0.108295654650418
$newSheet.Range("F67").Value = 0.064000000000000001
$newSheet.Range("A68").Value = "AD9_L.png"
$newSheet.Range("B68").Value = "CONTRALATERAL"
$newSheet.Range("C68").Value = 1
$newSheet.Range("D68").Value = 0
$newSheet.Range("E68").Value = 0.40783190376601097
$newSheet.Range("F68").Value = 0.43166666599999998

# --- Header row last; G1 ("Naive combined output") before D1
#     ("Radiologist output(89%)") to match new shared-string order ---
$newSheet.Range("A1").Value = "Image"
$newSheet.Range("B1").Value = "Ground truth"
$newSheet.Range("C1").Value = "Classifier output (66%)"
$newSheet.Range("E1").Value = "Classifier confidence"
$newSheet.Range("F1").Value = "Radiologist confidence"
$newSheet.Range("G1").Value = "Naïve combined output"
$newSheet.Range("D1").Value = "Radiologist output(89%)"

# --- Column widths (match autofit widths from the authored workbook) ---
$newSheet.Columns("A:A").ColumnWidth = 11.5703125
$newSheet.Columns("B:B").ColumnWidth = 16
$newSheet.Columns("C:C").ColumnWidth = 21.42578125
$newSheet.Columns("D:D").ColumnWidth = 22.7109375
$newSheet.Columns("E:E").ColumnWidth = 19.7109375
$newSheet.Columns("F:F").ColumnWidth = 21.5703125
$newSheet.Columns("G:G").ColumnWidth = 22.28515625

# --- Header row formatting (bold white text, grey fill, double border; same "Check Cell" style used by the other result sheets) ---
$newSheet.Range("A1:G1").Style = "Check Cell"

# --- View state: scrolled to row 19, K47 selected, this sheet active/tabSelected ---
$newSheet.Activate()
$excel.ActiveWindow.ScrollRow = 19
$newSheet.Range("K47").Select()

# --- "Better split" sheet view tidy-up: drop the stale selection/scroll, select row 1 ---
$ws7 = $wb.Worksheets.Item("Better split")
$ws7.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws7.Rows("1:1").Select()

# --- Re-activate the new sheet so it ends up tabSelected/active like the target workbook ---
$newSheet.Activate()
